# "Generate Report for Handoff"
#
# The localization-status report has moved on from "In Translation" to
# "Ready for handoff": refresh the status text and the two timestamps
# that were recorded at hand-off time, and widen the "Status" columns
# so the new (longer) status text fits.

$wb = $excel.ActiveWorkbook

$statusNew = "Ready for handoff"

# Column width (in "characters") that best approximates the new, wider
# column used to fit "Ready for handoff". ColumnWidth is quantized by
# Excel to a pixel grid, so this is the closest achievable value to the
# target stored width of ~17.216 characters.
$newStatusColWidth = 16.333

# ---- Overview sheet ----------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $statusNew
$wsOverview.Range("F2").Value = $statusNew
$wsOverview.Range("G2").Value = "2016-08-26 15:12:17"

$wsOverview.Columns.Item(5).ColumnWidth = $newStatusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusColWidth

# ---- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = $statusNew
$wsZhCn.Range("H2").Value = "2016-08-26 15:12:11"

$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusColWidth

# ---- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = $statusNew
$wsDeDe.Range("H2").Value = "2016-08-26 15:12:17"

$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusColWidth
